$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "_old" and "_new" header-column suffixes to "_FV2210" / "_FV2304"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2210")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2304")
}

# 2. Turn the data range into a real Excel Table (ListObject) with an AutoFilter
$rng = $ws.Range("A1:U70")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# 3. Freeze the header row (split/pane at row 2)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$null = $excel.ActiveWindow.FreezePanes
